$d = $word.ActiveDocument

# Position at the very end of the document (collapsed range after the
# last paragraph mark) and insert the new "4/4/23" meeting-notes block
# as raw WordprocessingML via InsertXML. Using InsertXML (rather than a
# sequence of InsertParagraphAfter/typing calls) lets us precisely
# reproduce the target markup, including reusing the existing bullet
# list definition (numId 2) instead of Word auto-creating a brand new
# list the way ApplyBulletDefault/ApplyListTemplateWithLevel would.
#
# Note: the document's very last paragraph mark can never be "pushed
# down" by an insertion collapsed at Content.End - the final <w:p> of
# whatever we insert there always fuses with the pre-existing last
# paragraph instead of creating an additional break. To still end up
# with two brand-new empty paragraphs (in addition to the two empty
# paragraphs already at the end of the document) before the dated
# block, we insert three leading empty paragraphs here.
$end = $d.Content
$end.Collapse(0)

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
<w:p/>
<w:p/>
<w:p><w:r><w:t>4/4/23</w:t></w:r></w:p>
<w:p><w:r><w:t>Present: Obaidullah Sarsour, Connor Lim, Phat Le</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Activities</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Discussed submission of Sprint 0 documents</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Planned initial steps of the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>project</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Discussed retrospective.</w:t></w:r></w:p>
<w:p><w:r><w:t>T</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>o do</w:t></w:r><w:r><w:t xml:space="preserve">:  </w:t></w:r></w:p>
<w:p><w:r><w:t>Come up with project ideas.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$end.InsertXML($xml)
